# Folienmaster hat jetzt das richtige Datum
# Update every "dt" (date) placeholder -- on the slide master, on every
# slide layout belonging to it, and on every slide -- from the stale
# "14.07.2016" to the corrected "10.07.2018".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Placeholders.Count; $i++) {
        $sh = $shapes.Placeholders.Item($i)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "10.07.2018"
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Every slide in the deck
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    Update-DatePlaceholder $slide.Shapes
}
